$d = $word.ActiveDocument

# --- Locate the paragraph that asks "How old were you when you learned
#     English? __________" and split it into two questions: one about
#     starting to learn English, one about feeling reasonably fluent.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "How old were you when you learned English?*") {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    # Create the second question paragraph right after the first one,
    # inheriting the same paragraph formatting (style/numbering/indent).
    $target.Range.InsertParagraphAfter()
    $second = $d.Paragraphs.Item($target.Index + 1)

    # Rewrite the first paragraph's wording.
    $findRange = $target.Range.Duplicate
    $findRange.Find.Execute("learned", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $findRange.Text = "started to learn"

    # Fill in the text of the newly inserted second paragraph.
    $secondStart = $second.Range.Start
    $insertPoint = $d.Range($secondStart, $secondStart)
    $insertPoint.InsertAfter("How old were you when you felt reasonably fluent in English? __________")
}
